$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows (2-9) and add new data rows (10-17) per the diff.
# Columns A-D are strings (shared-string backed); columns E-T are numeric.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt2b"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.921795333333333
$ws.Range("H2").Value = 5.765385999999999
$ws.Range("I2").Value = 0.1981502018195672
$ws.Range("J2").Value = 0.1981502018195673
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 19.48350866666667
$ws.Range("N2").Value = 58.450526
$ws.Range("O2").Value = 0.3081250754721727
$ws.Range("P2").Value = 0.3081250754721726
$ws.Range("Q2").Value = 37.44331603255955
$ws.Range("R2").Value = 336.989844293036
$ws.Range("S2").Value = 0.06105504589048039
$ws.Range("T2").Value = 0.0610550458904804

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt2b"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.921795333333333
$ws.Range("H3").Value = 5.765385999999999
$ws.Range("I3").Value = 0.1981502018195672
$ws.Range("J3").Value = 0.1981502018195673
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.343383
$ws.Range("N3").Value = 61.03014900000001
$ws.Range("O3").Value = 0.3217236961512193
$ws.Range("P3").Value = 0.3217236961512193
$ws.Range("Q3").Value = 39.09581851361267
$ws.Range("R3").Value = 351.862366622514
$ws.Range("S3").Value = 0.06374961532250123
$ws.Range("T3").Value = 0.06374961532250124

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt2b"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.921795333333333
$ws.Range("H4").Value = 5.765385999999999
$ws.Range("I4").Value = 0.1981502018195672
$ws.Range("J4").Value = 0.1981502018195673
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1305583333333333
$ws.Range("N4").Value = 0.391675
$ws.Range("O4").Value = 0.002064735720865253
$ws.Range("P4").Value = 0.002064735720865253
$ws.Range("Q4").Value = 0.2509063957277777
$ws.Range("R4").Value = 2.25815756155
$ws.Range("S4").Value = 0.0004091277997935195
$ws.Range("T4").Value = 0.0004091277997935196

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt2b"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.921795333333333
$ws.Range("H5").Value = 5.765385999999999
$ws.Range("I5").Value = 0.1981502018195672
$ws.Range("J5").Value = 0.1981502018195673
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 23.275017
$ws.Range("N5").Value = 69.825051
$ws.Range("O5").Value = 0.3680864926557428
$ws.Range("P5").Value = 0.3680864926557428
$ws.Range("Q5").Value = 44.72981905385399
$ws.Range("R5").Value = 402.568371484686
$ws.Range("S5").Value = 0.07293641280679208
$ws.Range("T5").Value = 0.0729364128067921

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt2b"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.24427
$ws.Range("H6").Value = 12.73281
$ws.Range("I6").Value = 0.437613174769253
$ws.Range("J6").Value = 0.437613174769253
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 19.48350866666667
$ws.Range("N6").Value = 58.450526
$ws.Range("O6").Value = 0.3081250754721727
$ws.Range("P6").Value = 0.3081250754721726
$ws.Range("Q6").Value = 82.69327132867333
$ws.Range("R6").Value = 744.23944195806
$ws.Range("S6").Value = 0.1348395925033932
$ws.Range("T6").Value = 0.1348395925033931

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt2b"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.24427
$ws.Range("H7").Value = 12.73281
$ws.Range("I7").Value = 0.437613174769253
$ws.Range("J7").Value = 0.437613174769253
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 20.343383
$ws.Range("N7").Value = 61.03014900000001
$ws.Range("O7").Value = 0.3217236961512193
$ws.Range("P7").Value = 0.3217236961512193
$ws.Range("Q7").Value = 86.34281016541001
$ws.Range("R7").Value = 777.0852914886901
$ws.Range("S7").Value = 0.1407905280712336
$ws.Range("T7").Value = 0.1407905280712336

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt2b"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.24427
$ws.Range("H8").Value = 12.73281
$ws.Range("I8").Value = 0.437613174769253
$ws.Range("J8").Value = 0.437613174769253
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1305583333333333
$ws.Range("N8").Value = 0.391675
$ws.Range("O8").Value = 0.002064735720865253
$ws.Range("P8").Value = 0.002064735720865253
$ws.Range("Q8").Value = 0.5541248174166666
$ws.Range("R8").Value = 4.987123356750001
$ws.Range("S8").Value = 0.0009035555538673256
$ws.Range("T8").Value = 0.0009035555538673256

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt2b"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.24427
$ws.Range("H9").Value = 12.73281
$ws.Range("I9").Value = 0.437613174769253
$ws.Range("J9").Value = 0.437613174769253
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 23.275017
$ws.Range("N9").Value = 69.825051
$ws.Range("O9").Value = 0.3680864926557428
$ws.Range("P9").Value = 0.3680864926557428
$ws.Range("Q9").Value = 98.78545640259001
$ws.Range("R9").Value = 889.0691076233101
$ws.Range("S9").Value = 0.1610794986407589
$ws.Range("T9").Value = 0.1610794986407589

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Wnt2b"
$ws.Range("C10").Value = "Fzd4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.068596333333333
$ws.Range("H10").Value = 3.205789
$ws.Range("I10").Value = 0.1101795677411623
$ws.Range("J10").Value = 0.1101795677411623
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 19.48350866666667
$ws.Range("N10").Value = 58.450526
$ws.Range("O10").Value = 0.3081250754721727
$ws.Range("P10").Value = 0.3081250754721726
$ws.Range("Q10").Value = 20.82000592166822
$ws.Range("R10").Value = 187.380053295014
$ws.Range("S10").Value = 0.03394908762573699
$ws.Range("T10").Value = 0.03394908762573699

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Wnt2b"
$ws.Range("C11").Value = "Fzd4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.068596333333333
$ws.Range("H11").Value = 3.205789
$ws.Range("I11").Value = 0.1101795677411623
$ws.Range("J11").Value = 0.1101795677411623
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 20.343383
$ws.Range("N11").Value = 61.03014900000001
$ws.Range("O11").Value = 0.3217236961512193
$ws.Range("P11").Value = 0.3217236961512193
$ws.Range("Q11").Value = 21.73886448139567
$ws.Range("R11").Value = 195.6497803325611
$ws.Range("S11").Value = 0.03544737777403038
$ws.Range("T11").Value = 0.03544737777403038

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Wnt2b"
$ws.Range("C12").Value = "Fzd4"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.068596333333333
$ws.Range("H12").Value = 3.205789
$ws.Range("I12").Value = 0.1101795677411623
$ws.Range("J12").Value = 0.1101795677411623
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1305583333333333
$ws.Range("N12").Value = 0.391675
$ws.Range("O12").Value = 0.002064735720865253
$ws.Range("P12").Value = 0.002064735720865253
$ws.Range("Q12").Value = 0.1395141562861111
$ws.Range("R12").Value = 1.255627406575
$ws.Range("S12").Value = 0.0002274916892246707
$ws.Range("T12").Value = 0.0002274916892246707

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Wnt2b"
$ws.Range("C13").Value = "Fzd4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.068596333333333
$ws.Range("H13").Value = 3.205789
$ws.Range("I13").Value = 0.1101795677411623
$ws.Range("J13").Value = 0.1101795677411623
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 23.275017
$ws.Range("N13").Value = 69.825051
$ws.Range("O13").Value = 0.3680864926557428
$ws.Range("P13").Value = 0.3680864926557428
$ws.Range("Q13").Value = 24.87159782447101
$ws.Range("R13").Value = 223.844380420239
$ws.Range("S13").Value = 0.04055561065217026
$ws.Range("T13").Value = 0.04055561065217025

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Wnt2b"
$ws.Range("C14").Value = "Fzd4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.464018
$ws.Range("H14").Value = 7.392054
$ws.Range("I14").Value = 0.2540570556700175
$ws.Range("J14").Value = 0.2540570556700175
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 19.48350866666667
$ws.Range("N14").Value = 58.450526
$ws.Range("O14").Value = 0.3081250754721727
$ws.Range("P14").Value = 0.3081250754721726
$ws.Range("Q14").Value = 48.00771605782266
$ws.Range("R14").Value = 432.069444520404
$ws.Range("S14").Value = 0.0782813494525621
$ws.Range("T14").Value = 0.0782813494525621

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Wnt2b"
$ws.Range("C15").Value = "Fzd4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.464018
$ws.Range("H15").Value = 7.392054
$ws.Range("I15").Value = 0.2540570556700175
$ws.Range("J15").Value = 0.2540570556700175
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 20.343383
$ws.Range("N15").Value = 61.03014900000001
$ws.Range("O15").Value = 0.3217236961512193
$ws.Range("P15").Value = 0.3217236961512193
$ws.Range("Q15").Value = 50.12646189289401
$ws.Range("R15").Value = 451.138157036046
$ws.Range("S15").Value = 0.0817361749834541
$ws.Range("T15").Value = 0.08173617498345413

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Wnt2b"
$ws.Range("C16").Value = "Fzd4"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.464018
$ws.Range("H16").Value = 7.392054
$ws.Range("I16").Value = 0.2540570556700175
$ws.Range("J16").Value = 0.2540570556700175
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.1305583333333333
$ws.Range("N16").Value = 0.391675
$ws.Range("O16").Value = 0.002064735720865253
$ws.Range("P16").Value = 0.002064735720865253
$ws.Range("Q16").Value = 0.3216980833833333
$ws.Range("R16").Value = 2.89528275045
$ws.Range("S16").Value = 0.0005245606779797372
$ws.Range("T16").Value = 0.0005245606779797373

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Wnt2b"
$ws.Range("C17").Value = "Fzd4"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.464018
$ws.Range("H17").Value = 7.392054
$ws.Range("I17").Value = 0.2540570556700175
$ws.Range("J17").Value = 0.2540570556700175
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 23.275017
$ws.Range("N17").Value = 69.825051
$ws.Range("O17").Value = 0.3680864926557428
$ws.Range("P17").Value = 0.3680864926557428
$ws.Range("Q17").Value = 57.350060838306
$ws.Range("R17").Value = 516.150547544754
$ws.Range("S17").Value = 0.09351497055602154
$ws.Range("T17").Value = 0.09351497055602154

Write-Host "Applied Wnt2b-Fzd4 NATMI update (rows 2-17)"